# MeetingRoomData.xlsx — "Reading and populating the meetingRoomList done"
#
# 1. Remove the obsolete "15-04-2018" sheet.
# 2. Remove the now-unused "Booked On" column (column B) from the two
#    remaining sheets (Template, 26-04-2018) — it was never actually used.
# 3. Populate the first few booked slots for "26-04-2018" with data
#    (room 1-10, slots 1-3 => 0, i.e. not booked).
# 4. Leave the workbook with "26-04-2018" as the active sheet/tab, with
#    the Template sheet scrolled/selected across its full data range and
#    26-04-2018 positioned at the next empty row underneath the table.

$wb = $excel.ActiveWorkbook

$wsTemplate = $wb.Worksheets.Item("Template")
$ws2604     = $wb.Worksheets.Item("26-04-2018")
$ws1504     = $wb.Worksheets.Item("15-04-2018")

# --- 1. Drop the old 15-04-2018 sheet -------------------------------------
$ws1504.Delete() | Out-Null

# --- 2. Drop the unused "Booked On" column on the remaining sheets --------
$wsTemplate.Columns("B").Delete() | Out-Null
$ws2604.Columns("B").Delete() | Out-Null

# --- 3. Fill in the first rooms/slots on the 26-04-2018 sheet -------------
$ws2604.Range("B2:D11").Value = 0

# --- 4. View state: which sheet/cells are selected ------------------------
$wsTemplate.Range("A1:AW11").Select() | Out-Null
$ws2604.Range("K16").Select() | Out-Null
$ws2604.Activate() | Out-Null
